# Add 2022-Q1 data
# 1) Insert a new worksheet "2022-Q1" right before the "总计" (totals) sheet,
#    populated with the quarter's fund-holding detail (same shape as the
#    other quarterly sheets).
# 2) Insert a new top data row into "总计" summarising the new quarter and
#    push the existing history rows down.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. New "2022-Q1" sheet
# ---------------------------------------------------------------------
$template = $wb.Worksheets.Item("2021-Q4")
$totalSheetRef = $wb.Worksheets.Item("总计")

$newSheet = $wb.Worksheets.Add($totalSheetRef)
$newSheet.Name = "2022-Q1"

# NOTE: inserting a sheet shifts worksheet positions, and the engine
# resolves already-bound sheet variables by position rather than fixed
# identity -- so re-fetch the "总计" sheet by name now that it has moved.
$totalSheet = $wb.Worksheets.Item("总计")

# Pull over the header / index-column formatting from the previous quarter
# sheet so the new sheet matches the established look (bold, centered,
# bordered header row + first data column).
$template.Range("B1:H1").Copy()
$newSheet.Range("B1:H1").PasteSpecial(-4122)
$template.Range("A2:A3").Copy()
$newSheet.Range("A2:A3").PasteSpecial(-4122)

$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

# The source data keeps these figures as plain text (not numbers), so
# format the cells as Text before writing them -- otherwise Excel would
# normalise "007107" to 7107 and "0.0690" to 0.069.
$newSheet.Range("B2:G2").NumberFormat = "@"
$newSheet.Range("B3:F3").NumberFormat = "@"

$newSheet.Range("A2").Value = 0
$newSheet.Range("B2").Value = "007107"
$newSheet.Range("C2").Value = "太平 MSCI 香港价值增强指数A"
$newSheet.Range("D2").Value = "1.05"
$newSheet.Range("E2").Value = "93.78"
$newSheet.Range("F2").Value = "6.57"
$newSheet.Range("G2").Value = "0.0690"
$newSheet.Range("H2").Value = 3

$newSheet.Range("A3").Value = 1
$newSheet.Range("B3").Value = "007108"
$newSheet.Range("C3").Value = "太平 MSCI 香港价值增强指数C"
$newSheet.Range("D3").Value = "0.00"
$newSheet.Range("E3").Value = "93.78"
$newSheet.Range("F3").Value = "6.57"
$newSheet.Range("G3").Value = 0
$newSheet.Range("H3").Value = 3

# ---------------------------------------------------------------------
# 2. Add the new quarter to the "总计" summary sheet
# ---------------------------------------------------------------------
$totalSheet.Rows.Item(2).Insert()
$totalSheet.Range("A2:D2").ClearFormats()

$totalSheet.Range("A3").Copy()
$totalSheet.Range("A2").PasteSpecial(-4122)

$totalSheet.Range("A2").Value = 0
$totalSheet.Range("B2").Value = "2022-Q1"
$totalSheet.Range("C2").Value = 2
$totalSheet.Range("D2").Value = 0.07
